$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 766318
$ws.Range("E2").Value = 1429203416

# Row 69
$ws.Range("C69").Value = 17890
$ws.Range("E69").Value = 103828651

# Row 88
$ws.Range("C88").Value = 71270
$ws.Range("E88").Value = 110301078

# Row 91
$ws.Range("C91").Value = 18866
$ws.Range("E91").Value = 75241295

# Row 112
$ws.Range("C112").Value = 145227
$ws.Range("E112").Value = 716160483

# Row 121
$ws.Range("C121").Value = 1306206
$ws.Range("E121").Value = 2274880771

# Row 129
$ws.Range("C129").Value = 633508
$ws.Range("E129").Value = 3429922615

# Row 130
$ws.Range("C130").Value = 4241
$ws.Range("E130").Value = 140443816

# Row 132
$ws.Range("C132").Value = 585779
$ws.Range("E132").Value = 3465773566

# Row 136
$ws.Range("C136").Value = 26686
$ws.Range("D136").Value = 4273
$ws.Range("E136").Value = 143736140

# Row 139
$ws.Range("C139").Value = 76643
$ws.Range("E139").Value = 114133158

# Row 151
$ws.Range("C151").Value = 39931
$ws.Range("E151").Value = 60378307

# Row 154
$ws.Range("C154").Value = 18459
$ws.Range("E154").Value = 73371235

# Row 178
$ws.Range("C178").Value = 515878
$ws.Range("E178").Value = 891190220

# Row 221
$ws.Range("C221").Value = 135497
$ws.Range("D221").Value = 27176
$ws.Range("E221").Value = 681872566
